$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 8.88020426756265
$ws.Range("D2").Value = 7.31260311002256
$ws.Range("E2").Value = 13.17366153915612
$ws.Range("F2").Value = 42.72123245259331
$ws.Range("G2").Value = 51.14327410568792
$ws.Range("H2").Value = 19.85655135252868
$ws.Range("I2").Value = 32.95461377668224
$ws.Range("J2").Value = 10.61427242236585
$ws.Range("K2").Value = 20.35528468791982
$ws.Range("L2").Value = 10.44912945570019
$ws.Range("N2").Value = 18.98221975678306

$ws.Range("C3").Value = 8.84619330212686
$ws.Range("D3").Value = 7.291746680852307
$ws.Range("E3").Value = 13.17509751968188
$ws.Range("F3").Value = 42.77946693131241
$ws.Range("G3").Value = 51.19074642049652
$ws.Range("H3").Value = 19.91338659552721
$ws.Range("I3").Value = 33.01745030940503
$ws.Range("J3").Value = 10.63554970605779
$ws.Range("K3").Value = 20.02089748042232
$ws.Range("L3").Value = 10.46643927751476
$ws.Range("N3").Value = 19.04134677326683

$ws.Range("C4").Value = 8.826535753964091
$ws.Range("D4").Value = 7.279619139429999
$ws.Range("E4").Value = 13.1777976426032
$ws.Range("F4").Value = 42.82585489569176
$ws.Range("G4").Value = 51.23536699543946
$ws.Range("H4").Value = 19.95213078219293
$ws.Range("I4").Value = 33.06367098964017
$ws.Range("J4").Value = 10.64977244910957
$ws.Range("K4").Value = 19.8162816067607
$ws.Range("L4").Value = 10.47815221363428
$ws.Range("N4").Value = 19.07961914397178

$ws.Range("C5").Value = 8.818838679814526
$ws.Range("D5").Value = 7.274850511660517
$ws.Range("E5").Value = 13.17935626254699
$ws.Range("F5").Value = 42.84742466629696
$ws.Range("G5").Value = 51.25742736973332
$ws.Range("H5").Value = 19.96888477032422
$ws.Range("I5").Value = 33.08442184746694
$ws.Range("J5").Value = 10.65585995109731
$ws.Range("K5").Value = 19.73317845945127
$ws.Range("L5").Value = 10.48319844050587
$ws.Range("N5").Value = 19.09571159889364

$ws.Range("C6").Value = 8.817579675378637
$ws.Range("D6").Value = 7.274069241263102
$ws.Range("E6").Value = 13.17964277805867
$ws.Range("F6").Value = 42.85116711949897
$ws.Range("G6").Value = 51.26132420746008
$ws.Range("H6").Value = 19.97172500724083
$ws.Range("I6").Value = 33.08798304231758
$ws.Range("J6").Value = 10.65688839880722
$ws.Range("K6").Value = 19.71939916682585
$ws.Range("L6").Value = 10.48405286853624
$ws.Range("N6").Value = 19.09841374220655

$ws.Range("C7").Value = 8.826430672102816
$ws.Range("D7").Value = 7.279554122089174
$ws.Range("E7").Value = 13.17781680578565
$ws.Range("F7").Value = 42.82613500770865
$ws.Range("G7").Value = 51.23564883070505
$ws.Range("H7").Value = 19.95235282609866
$ws.Range("I7").Value = 33.06394309437407
$ws.Range("J7").Value = 10.6498533661064
$ws.Range("K7").Value = 19.81515957964769
$ws.Range("L7").Value = 10.47821916251699
$ws.Range("N7").Value = 19.07983416189605

$ws.Range("C8").Value = 8.868226315546639
$ws.Range("D8").Value = 7.305272761935511
$ws.Range("E8").Value = 13.17377991133387
$ws.Range("F8").Value = 42.73910107790761
$ws.Range("G8").Value = 51.15642380299204
$ws.Range("H8").Value = 19.87534850818717
$ws.Range("I8").Value = 32.97469134816794
$ws.Range("J8").Value = 10.62136860102109
$ws.Range("K8").Value = 20.2399047823583
$ws.Range("L8").Value = 10.45487298999463
$ws.Range("N8").Value = 19.00219901199261

$ws.Range("C9").Value = 8.959658613855373
$ws.Range("D9").Value = 7.360959905831923
$ws.Range("E9").Value = 13.18024233598129
$ws.Range("F9").Value = 42.65306443582478
$ws.Range("G9").Value = 51.124334635881
$ws.Range("H9").Value = 19.75495476286142
$ws.Range("I9").Value = 32.86049652206117
$ws.Range("J9").Value = 10.57468700647669
$ws.Range("K9").Value = 21.07388015762417
$ws.Range("L9").Value = 10.41768102932987
$ws.Range("N9").Value = 18.86551664178042

$ws.Range("C10").Value = 9.03225317826822
$ws.Range("D10").Value = 7.404895773596105
$ws.Range("E10").Value = 13.19368352401292
$ws.Range("F10").Value = 42.64176935486302
$ws.Range("G10").Value = 51.17644798609883
$ws.Range("H10").Value = 19.68527837651138
$ws.Range("I10").Value = 32.81394801394606
$ws.Range("J10").Value = 10.54596499628297
$ws.Range("K10").Value = 21.68135649543964
$ws.Range("L10").Value = 10.39557072000331
$ws.Range("N10").Value = 18.77450362900549

$ws.Range("C11").Value = 9.066373878596517
$ws.Range("D11").Value = 7.425499084768152
$ws.Range("E11").Value = 13.20166769878634
$ws.Range("F11").Value = 42.64794898739201
$ws.Range("G11").Value = 51.21665383179423
$ws.Range("H11").Value = 19.65768131388248
$ws.Range("I11").Value = 32.80092811687874
$ws.Range("J11").Value = 10.53410498186516
$ws.Range("K11").Value = 21.95538213171602
$ws.Range("L11").Value = 10.38663969071774
$ws.Range("N11").Value = 18.73512594511838

$ws.Range("C12").Value = 9.079444996033995
$ws.Range("D12").Value = 7.433386072976951
$ws.Range("E12").Value = 13.20495813607705
$ws.Range("F12").Value = 42.65191809309419
$ws.Range("G12").Value = 51.23425194995955
$ws.Range("H12").Value = 19.64782225362009
$ws.Range("I12").Value = 32.79717322104121
$ws.Range("J12").Value = 10.52978695606281
$ws.Range("K12").Value = 22.05872048661547
$ws.Range("L12").Value = 10.38341940714562
$ws.Range("N12").Value = 18.72050459427998

$ws.Range("C13").Value = 9.076623330451669
$ws.Range("D13").Value = 7.431683750413592
$ws.Range("E13").Value = 13.20423763764747
$ws.Range("F13").Value = 42.65099081349135
$ws.Range("G13").Value = 51.2303563605056
$ws.Range("H13").Value = 19.64991924889706
$ws.Range("I13").Value = 32.79792958951214
$ws.Range("J13").Value = 10.53070922588647
$ws.Range("K13").Value = 22.03648529255254
$ws.Range("L13").Value = 10.3841057669676
$ws.Range("N13").Value = 18.72364067625636

$ws.Range("C14").Value = 9.067446277275591
$ws.Range("D14").Value = 7.42614627089073
$ws.Range("E14").Value = 13.20193306506393
$ws.Range("F14").Value = 42.64824287601177
$ws.Range("G14").Value = 51.21805408741634
$ws.Range("H14").Value = 19.65685834097326
$ws.Range("I14").Value = 32.80059562803626
$ws.Range("J14").Value = 10.53374626765799
$ws.Range("K14").Value = 21.96389297610681
$ws.Range("L14").Value = 10.38637151740663
$ws.Range("N14").Value = 18.73391722896755

$ws.Range("C15").Value = 9.061844423528846
$ws.Range("D15").Value = 7.422765357426559
$ws.Range("E15").Value = 13.20055616187662
$ws.Range("F15").Value = 42.64677185427689
$ws.Range("G15").Value = 51.21082759661222
$ws.Range("H15").Value = 19.66118580123386
$ws.Range("I15").Value = 32.80238180309032
$ws.Range("J15").Value = 10.53562907821989
$ws.Range("K15").Value = 21.91936941905142
$ws.Range("L15").Value = 10.38778040248671
$ws.Range("N15").Value = 18.74024966736088

$ws.Range("C16").Value = 9.030044757792487
$ws.Range("D16").Value = 7.403561392322748
$ws.Range("E16").Value = 13.19319920728872
$ws.Range("F16").Value = 42.64159340379679
$ws.Range("G16").Value = 51.17415254852309
$ws.Range("H16").Value = 19.68716455840244
$ws.Range("I16").Value = 32.81496324838368
$ws.Range("J16").Value = 10.54676431511516
$ws.Range("K16").Value = 21.66339333697524
$ws.Range("L16").Value = 10.39617703199512
$ws.Range("N16").Value = 18.77711770620003

$ws.Range("C17").Value = 9.010812477837549
$ws.Range("D17").Value = 7.391935712541459
$ws.Range("E17").Value = 13.18916355648063
$ws.Range("F17").Value = 42.64131684019886
$ws.Range("G17").Value = 51.15588064129928
$ws.Range("H17").Value = 19.70415294163739
$ws.Range("I17").Value = 32.82477246750981
$ws.Range("J17").Value = 10.55390404831021
$ws.Range("K17").Value = 21.50569701387284
$ws.Range("L17").Value = 10.40161650172751
$ws.Range("N17").Value = 18.80025287257779

$ws.Range("C18").Value = 8.999854298448456
$ws.Range("D18").Value = 7.385307201029542
$ws.Range("E18").Value = 13.18701846355201
$ws.Range("F18").Value = 42.64222304270816
$ws.Range("G18").Value = 51.1469242527116
$ws.Range("H18").Value = 19.71430992623148
$ws.Range("I18").Value = 32.83118195801661
$ws.Range("J18").Value = 10.55812414419054
$ws.Range("K18").Value = 21.41478152058499
$ws.Range("L18").Value = 10.40485124769317
$ws.Range("N18").Value = 18.8137502310969

$ws.Range("C19").Value = 8.996162083203732
$ws.Range("D19").Value = 7.38307302077299
$ws.Range("E19").Value = 13.18632247263238
$ws.Range("F19").Value = 42.64271276236556
$ws.Range("G19").Value = 51.14415846597841
$ws.Range("H19").Value = 19.71781509422667
$ws.Range("I19").Value = 32.83348382406255
$ws.Range("J19").Value = 10.55957250023639
$ws.Range("K19").Value = 21.3839655463837
$ws.Range("L19").Value = 10.40596471207284
$ws.Range("N19").Value = 18.81835297569342

$ws.Range("C20").Value = 9.012849105290147
$ws.Range("D20").Value = 7.393167281676708
$ws.Range("E20").Value = 13.18957494501634
$ws.Range("F20").Value = 42.64123601383011
$ws.Range("G20").Value = 51.15766496178578
$ws.Range("H20").Value = 19.70230456494342
$ws.Range("I20").Value = 32.82364880538794
$ws.Range("J20").Value = 10.55313226568939
$ws.Range("K20").Value = 21.52250678177176
$ws.Range("L20").Value = 10.40102648185366
$ws.Range("N20").Value = 18.79777037326991

$ws.Range("C21").Value = 9.070137779928903
$ws.Range("D21").Value = 7.427770487097604
$ws.Range("E21").Value = 13.20260274304061
$ws.Range("F21").Value = 42.64900579500841
$ws.Range("G21").Value = 51.22160316951226
$ws.Range("H21").Value = 19.65480409766577
$ws.Range("I21").Value = 32.7997806276953
$ws.Range("J21").Value = 10.53284951915467
$ws.Range("K21").Value = 21.98522749030118
$ws.Range("L21").Value = 10.38570162596459
$ws.Range("N21").Value = 18.73089089037278

$ws.Range("C22").Value = 9.108452006083914
$ws.Range("D22").Value = 7.450879114509784
$ws.Range("E22").Value = 13.21267249696256
$ws.Range("F22").Value = 42.66357874485316
$ws.Range("G22").Value = 51.27722064244025
$ws.Range("H22").Value = 19.62720744509931
$ws.Range("I22").Value = 32.79103366681081
$ws.Range("J22").Value = 10.52060242134369
$ws.Range("K22").Value = 22.28509908604212
$ws.Range("L22").Value = 10.37662830109578
$ws.Range("N22").Value = 18.68887189483129

$ws.Range("C23").Value = 9.087925625890076
$ws.Range("D23").Value = 7.438501707105201
$ws.Range("E23").Value = 13.20715641873148
$ws.Range("F23").Value = 42.6549319216471
$ws.Range("G23").Value = 51.24627164301461
$ws.Range("H23").Value = 19.64162023244432
$ws.Range("I23").Value = 32.79507436317675
$ws.Range("J23").Value = 10.52704671302957
$ws.Range("K23").Value = 22.12531521662564
$ws.Range("L23").Value = 10.38138480341036
$ws.Range("N23").Value = 18.71114385890618

$ws.Range("C24").Value = 9.011928037996981
$ws.Range("D24").Value = 7.392610316923361
$ws.Range("E24").Value = 13.18938841084682
$ws.Range("F24").Value = 42.64126923739249
$ws.Range("G24").Value = 51.15685344731128
$ws.Range("H24").Value = 19.70313900086052
$ws.Range("I24").Value = 32.82415441452295
$ws.Range("J24").Value = 10.55348082922226
$ws.Range("K24").Value = 21.51490786876678
$ws.Range("L24").Value = 10.40129289488578
$ws.Range("N24").Value = 18.7988920985546

$ws.Range("C25").Value = 8.933945269566818
$ws.Range("D25").Value = 7.345350787650267
$ws.Range("E25").Value = 13.17696163884362
$ws.Range("F25").Value = 42.66724121833727
$ws.Range("G25").Value = 51.11975287822066
$ws.Range("H25").Value = 19.78423508415262
$ws.Range("I25").Value = 32.88484728656489
$ws.Range("J25").Value = 10.57468700647669
$ws.Range("K25").Value = 21.07388015762417
$ws.Range("L25").Value = 10.41768102932987
$ws.Range("N25").Value = 18.90083523766111
